# Auto-generated edit script applying numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns across
# several sheets, per scheduled runner update.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 749.5
$ws.Range("I4").Value = 749.5
$ws.Range("K4").Value = 749.5
$ws.Range("M4").Value = -635.5
$ws.Range("H17").Value = 1317.8667
$ws.Range("J17").Value = 1096.3864
$ws.Range("L17").Value = 3289.1592
$ws.Range("N17").Value = -3625.1592
$ws.Range("H43").Value = 1532.7693
$ws.Range("I43").Value = 1280
$ws.Range("J43").Value = 1553.8334
$ws.Range("K43").Value = 1280
$ws.Range("L43").Value = 1553.8334
$ws.Range("M43").Value = -1211
$ws.Range("N43").Value = -1691.8334
$ws.Range("H64").Value = 3500
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3500
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("H70").Value = 34000
$ws.Range("I70").Value = 500
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1230
$ws.Range("H73").Value = 34000
$ws.Range("I73").Value = 500
$ws.Range("K73").Value = 1500
$ws.Range("M73").Value = -564
$ws.Range("H74").Value = 5402.3335
$ws.Range("I74").Value = 5745.857
$ws.Range("K74").Value = 5745.857
$ws.Range("M74").Value = -4809.857
$ws.Range("H77").Value = 5402.3335
$ws.Range("I77").Value = 5745.857
$ws.Range("K77").Value = 28729.285
$ws.Range("M77").Value = -24049.285
$ws.Range("H94").Value = 3045.6
$ws.Range("I94").Value = 2557
$ws.Range("K94").Value = 2557
$ws.Range("M94").Value = -2106
$ws.Range("H106").Value = 4048.1538
$ws.Range("I106").Value = 2851.75
$ws.Range("K106").Value = 2851.75
$ws.Range("M106").Value = -2220.75
# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 28341
$ws.Range("I23").Value = 70006
$ws.Range("J23").Value = 17924.75
$ws.Range("K23").Value = 70006
$ws.Range("L23").Value = 17924.75
$ws.Range("M23").Value = -69747
$ws.Range("N23").Value = -18442.75
$ws.Range("H32").Value = 4754.4316
$ws.Range("I32").Value = 2771.5151
$ws.Range("J32").Value = 10703.182
$ws.Range("K32").Value = 2771.5151
$ws.Range("L32").Value = 10703.182
$ws.Range("M32").Value = -2484.5151
$ws.Range("N32").Value = -11277.182
$ws.Range("H122").Value = 2582.5
$ws.Range("I122").Value = 2633.3333
$ws.Range("J122").Value = 2430
$ws.Range("K122").Value = 7899.999899999999
$ws.Range("L122").Value = 7290
$ws.Range("M122").Value = -5449.999899999999
$ws.Range("N122").Value = -12190
$ws.Range("H132").Value = 1290.1538
$ws.Range("I132").Value = 1123.6957
$ws.Range("J132").Value = 2566.3333
$ws.Range("K132").Value = 3371.0871
$ws.Range("L132").Value = 7698.999899999999
$ws.Range("M132").Value = -841.0870999999997
$ws.Range("N132").Value = -12758.9999
# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2296.3635
$ws.Range("I20").Value = 2194.9
$ws.Range("J20").Value = 3311
$ws.Range("K20").Value = 2194.9
$ws.Range("L20").Value = 3311
$ws.Range("M20").Value = -1947.9
$ws.Range("N20").Value = -3805
$ws.Range("H86").Value = 112666
$ws.Range("I86").Value = 1377
$ws.Range("K86").Value = 1377
$ws.Range("M86").Value = -254
$ws.Range("H89").Value = 112666
$ws.Range("I89").Value = 1377
$ws.Range("K89").Value = 6885
$ws.Range("M89").Value = -1269
$ws.Range("H105").Value = 1941.381
$ws.Range("I105").Value = 2024.6842
$ws.Range("J105").Value = 1150
$ws.Range("K105").Value = 2024.6842
$ws.Range("L105").Value = 1150
$ws.Range("M105").Value = -277.6841999999999
$ws.Range("N105").Value = -4644
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2920.2222
$ws.Range("I31").Value = 1389.3043
$ws.Range("K31").Value = 1389.3043
$ws.Range("M31").Value = -1094.3043
$ws.Range("H32").Value = 2195.8
$ws.Range("I32").Value = 1494.75
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1494.75
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -1178.75
$ws.Range("N32").Value = -5632
$ws.Range("H34").Value = 2920.2222
$ws.Range("I34").Value = 1389.3043
$ws.Range("K34").Value = 1389.3043
$ws.Range("M34").Value = -1187.3043
$ws.Range("H45").Value = 4165.8335
$ws.Range("I45").Value = 4165.8335
$ws.Range("K45").Value = 4165.8335
$ws.Range("M45").Value = -3572.8335
$ws.Range("H86").Value = 58825270
$ws.Range("I86").Value = 90910570
$ws.Range("K86").Value = 90910570
$ws.Range("M86").Value = -90909447
$ws.Range("H89").Value = 58825270
$ws.Range("I89").Value = 90910570
$ws.Range("K89").Value = 454552850
$ws.Range("M89").Value = -454547234
$ws.Range("H94").Value = 1002.625
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1003
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1003
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -1905
$ws.Range("H107").Value = 343.5
$ws.Range("I107").Value = 343.5
$ws.Range("K107").Value = 343.5
$ws.Range("M107").Value = 1576.5
$ws.Range("H132").Value = 1732.6666
$ws.Range("I132").Value = 1084.875
$ws.Range("J132").Value = 3028.25
$ws.Range("K132").Value = 3254.625
$ws.Range("L132").Value = 9084.75
$ws.Range("M132").Value = -724.625
$ws.Range("N132").Value = -14144.75
$ws.Range("H134").Value = 2045.6875
$ws.Range("I134").Value = 1515.4667
$ws.Range("K134").Value = 4546.4001
$ws.Range("M134").Value = -2011.4001
# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 763.82355
$ws.Range("J107").Value = 763.82355
$ws.Range("L107").Value = 2291.47065
$ws.Range("N107").Value = -6131.470649999999
$ws.Range("H122").Value = 706.26666
$ws.Range("I122").Value = 504.85715
$ws.Range("J122").Value = 882.5
$ws.Range("K122").Value = 4543.71435
$ws.Range("L122").Value = 7942.5
$ws.Range("M122").Value = -2093.71435
$ws.Range("N122").Value = -12842.5
$ws.Range("H131").Value = 12422.55
$ws.Range("J131").Value = 13745.296
$ws.Range("L131").Value = 41235.888
$ws.Range("N131").Value = -51315.888
# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 15010
$ws.Range("J48").Value = 15010
$ws.Range("L48").Value = 15010
$ws.Range("N48").Value = -15980
$ws.Range("H80").Value = 2326.5
$ws.Range("I80").Value = 1995
$ws.Range("J80").Value = 2468.5715
$ws.Range("K80").Value = 1995
$ws.Range("L80").Value = 2468.5715
$ws.Range("M80").Value = -997
$ws.Range("N80").Value = -4464.5715
$ws.Range("H83").Value = 2326.5
$ws.Range("I83").Value = 1995
$ws.Range("J83").Value = 2468.5715
$ws.Range("K83").Value = 9975
$ws.Range("L83").Value = 12342.8575
$ws.Range("M83").Value = -4983
$ws.Range("N83").Value = -22326.8575
$ws.Range("H102").Value = 2853.9443
$ws.Range("I102").Value = 3651.1428
$ws.Range("K102").Value = 3651.1428
$ws.Range("M102").Value = -2029.1428
# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4124.5
$ws.Range("I22").Value = 4124.5
$ws.Range("K22").Value = 4124.5
$ws.Range("M22").Value = -3829.5
$ws.Range("H27").Value = 4124.5
$ws.Range("I27").Value = 4124.5
$ws.Range("K27").Value = 4124.5
$ws.Range("M27").Value = -4017.5
$ws.Range("H93").Value = 765.3333
$ws.Range("I93").Value = 427.8
$ws.Range("J93").Value = 1187.25
$ws.Range("K93").Value = 427.8
$ws.Range("L93").Value = 1187.25
$ws.Range("M93").Value = 820.2
$ws.Range("N93").Value = -3683.25
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459
$ws.Range("H110").Value = 26322
$ws.Range("J110").Value = 26322
$ws.Range("L110").Value = 26322
$ws.Range("N110").Value = -34502
# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11699.4
$ws.Range("J41").Value = 11699.4
$ws.Range("L41").Value = 11699.4
$ws.Range("N41").Value = -12479.4
$ws.Range("H100").Value = 486
$ws.Range("I100").Value = 329
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 658
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -117
$ws.Range("N100").Value = -2682
$ws.Range("H107").Value = 457.7647
$ws.Range("I107").Value = 352.13333
$ws.Range("K107").Value = 1056.39999
$ws.Range("M107").Value = 863.6000100000001
$ws.Range("H126").Value = 5357.0527
$ws.Range("I126").Value = 4653.231
$ws.Range("K126").Value = 13959.693
$ws.Range("M126").Value = -11489.693
$ws.Range("H132").Value = 7799.885
$ws.Range("I132").Value = 1959.8
$ws.Range("K132").Value = 5879.4
$ws.Range("M132").Value = -3349.4
